$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for the two additional columns (AC, AD)
$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

# New data values for row 2
$ws.Range("AC2").Value = 1
$ws.Range("AD2").Value = 1

# New data values for row 3
$ws.Range("AC3").Value = 1
$ws.Range("AD3").Value = 0.93856655290102387

# Extend the selection/used range to cover the newly added columns
$ws.Range("A1:AD3").Select() | Out-Null
